$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing shared-string values before they are overwritten,
# so we can reuse the same shared-string entries (no new unique strings).
$ECs   = $ws.Range("A2").Value()
$FAPs  = $ws.Range("A4").Value()
$MuSCs = $ws.Range("A6").Value()
$Fn1   = $ws.Range("B2").Value()
$Itgb8 = $ws.Range("C2").Value()

# Row 2: ECs -> ECs
$ws.Range("A2").Value = $ECs
$ws.Range("B2").Value = $Fn1
$ws.Range("C2").Value = $Itgb8
$ws.Range("D2").Value = $ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.027767
$ws.Range("N2").Value = 0.083301
$ws.Range("O2").Value = 0.002923627791763407
$ws.Range("P2").Value = 0.002923627791763407
$ws.Range("Q2").Value = 0.173269078836
$ws.Range("R2").Value = 1.559421709524
$ws.Range("S2").Value = 0.00005064397300360597
$ws.Range("T2").Value = 0.00005064397300360598

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = $ECs
$ws.Range("B3").Value = $Fn1
$ws.Range("C3").Value = $Itgb8
$ws.Range("D3").Value = $FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.237840333333334
$ws.Range("N3").Value = 12.713521
$ws.Range("O3").Value = 0.4462083687682946
$ws.Range("P3").Value = 0.4462083687682946
$ws.Range("Q3").Value = 26.444581366756
$ws.Range("R3").Value = 238.001232300804
$ws.Range("S3").Value = 0.007729357562391539
$ws.Range("T3").Value = 0.00772935756239154

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = $ECs
$ws.Range("B4").Value = $Fn1
$ws.Range("C4").Value = $Itgb8
$ws.Range("D4").Value = $MuSCs
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.231839666666667
$ws.Range("N4").Value = 15.695519
$ws.Range("O4").Value = 0.5508680034399419
$ws.Range("P4").Value = 0.5508680034399419
$ws.Range("Q4").Value = 32.64724455868399
$ws.Range("R4").Value = 293.825201028156
$ws.Range("S4").Value = 0.009542303699998613
$ws.Range("T4").Value = 0.009542303699998615

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = $FAPs
$ws.Range("B5").Value = $Fn1
$ws.Range("C5").Value = $Itgb8
$ws.Range("D5").Value = $ECs
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.027767
$ws.Range("N5").Value = 0.083301
$ws.Range("O5").Value = 0.002923627791763407
$ws.Range("P5").Value = 0.002923627791763407
$ws.Range("Q5").Value = 9.595347199093
$ws.Range("R5").Value = 86.358124791837
$ws.Range("S5").Value = 0.002804577180046318
$ws.Range("T5").Value = 0.002804577180046319

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = $FAPs
$ws.Range("B6").Value = $Fn1
$ws.Range("C6").Value = $Itgb8
$ws.Range("D6").Value = $FAPs
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.237840333333334
$ws.Range("N6").Value = 12.713521
$ws.Range("O6").Value = 0.4462083687682946
$ws.Range("P6").Value = 0.4462083687682946
$ws.Range("Q6").Value = 1464.45598633822
$ws.Range("R6").Value = 13180.10387704398
$ws.Range("S6").Value = 0.4280386895072046
$ws.Range("T6").Value = 0.4280386895072046

# Row 7: FAPs -> MuSCs
$ws.Range("A7").Value = $FAPs
$ws.Range("B7").Value = $Fn1
$ws.Range("C7").Value = $Itgb8
$ws.Range("D7").Value = $MuSCs
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.231839666666667
$ws.Range("N7").Value = 15.695519
$ws.Range("O7").Value = 0.5508680034399419
$ws.Range("P7").Value = 0.5508680034399419
$ws.Range("Q7").Value = 1807.9489354865
$ws.Range("R7").Value = 16271.5404193785
$ws.Range("S7").Value = 0.528436566384358
$ws.Range("T7").Value = 0.5284365663843581

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = $MuSCs
$ws.Range("B8").Value = $Fn1
$ws.Range("C8").Value = $Itgb8
$ws.Range("D8").Value = $ECs
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.027767
$ws.Range("N8").Value = 0.083301
$ws.Range("O8").Value = 0.002923627791763407
$ws.Range("P8").Value = 0.002923627791763407
$ws.Range("Q8").Value = 0.2340407865573333
$ws.Range("R8").Value = 2.106367079016
$ws.Range("S8").Value = 0.00006840663871348322
$ws.Range("T8").Value = 0.00006840663871348325

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = $MuSCs
$ws.Range("B9").Value = $Fn1
$ws.Range("C9").Value = $Itgb8
$ws.Range("D9").Value = $FAPs
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.237840333333334
$ws.Range("N9").Value = 12.713521
$ws.Range("O9").Value = 0.4462083687682946
$ws.Range("P9").Value = 0.4462083687682946
$ws.Range("Q9").Value = 35.71964868072622
$ws.Range("R9").Value = 321.476838126536
$ws.Range("S9").Value = 0.01044032169869848
$ws.Range("T9").Value = 0.01044032169869848

# Row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = $MuSCs
$ws.Range("B10").Value = $Fn1
$ws.Range("C10").Value = $Itgb8
$ws.Range("D10").Value = $MuSCs
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.231839666666667
$ws.Range("N10").Value = 15.695519
$ws.Range("O10").Value = 0.5508680034399419
$ws.Range("P10").Value = 0.5508680034399419
$ws.Range("Q10").Value = 44.09780929623378
$ws.Range("R10").Value = 396.880283666104
$ws.Range("S10").Value = 0.0128891333555853
$ws.Range("T10").Value = 0.01288913335558531
